# Update countries & provincias Spain
# Refresh COVID figures for the affected countries (the data refresh also
# re-sorts a handful of neighbouring rows by "Casos totales" descending,
# which is why some rows change which country they display) and bump the
# "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => updated country name + Casos totales / Nuevos casos / Casos activos /
# Recuperados / Casos criticos / Muertes hoy / Muertes
$updates = @(
    @{ Row = 8;   Name = "Alemania";     B = 147103; C = 38;  D = 95200; E = 47041; F = 2889; G = 0;   H = 4862 },
    @{ Row = 17;  Name = "Paises Bajos"; B = 34134;  C = 729; D = 0;     E = 29968; F = 1158; G = 165; H = 3916 },
    @{ Row = 19;  Name = "Portugal";     B = 21379;  C = 516; D = 917;   E = 19700; F = 213;  G = 27;  H = 762 },
    @{ Row = 20;  Name = "India";        B = 18985;  C = 446; D = 3273;  E = 15109; F = 0;    G = 11;  H = 603 },
    @{ Row = 23;  Name = "Suecia";       B = 15322;  C = 545; D = 550;   E = 13007; F = 515;  G = 185; H = 1765 },
    @{ Row = 24;  Name = "Austria";      B = 14873;  C = 78;  D = 10971; E = 3411;  F = 196;  G = 21;  H = 491 },
    @{ Row = 57;  Name = "Argentina";    B = 3031;   C = 0;   D = 840;   E = 2046;  F = 123;  G = 3;   H = 145 },
    @{ Row = 66;  Name = "Croacia";      B = 1908;   C = 27;  D = 801;   E = 1059;  F = 18;   G = 1;   H = 48 },
    @{ Row = 90;  Name = "Letonia";      B = 748;    C = 9;   D = 133;   E = 606;   F = 3;    G = 4;   H = 9 },
    @{ Row = 104; Name = "San Marino";   B = 476;    C = 14;  D = 62;    E = 374;   F = 4;    G = 1;   H = 40 },
    @{ Row = 109; Name = "Senegal";      B = 412;    C = 35;  D = 242;   E = 165;   F = 1;    G = 0;   H = 5 },
    @{ Row = 110; Name = "Georgia";      B = 408;    C = 6;   D = 95;    E = 309;   F = 6;    G = 0;   H = 4 },
    @{ Row = 111; Name = "Reunion";      B = 408;    C = 0;   D = 238;   E = 170;   F = 2;    G = 0;   H = 0 },
    @{ Row = 168; Name = "Siria";        B = 39;     C = 0;   D = 5;     E = 31;    F = 0;    G = 0;   H = 3 },
    @{ Row = 169; Name = "Mozambique";   B = 39;     C = 0;   D = 8;     E = 31;    F = 0;    G = 0;   H = 0 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.Name
    $ws.Range("B$r").Value = $u.B
    $ws.Range("C$r").Value = $u.C
    $ws.Range("D$r").Value = $u.D
    $ws.Range("E$r").Value = $u.E
    $ws.Range("F$r").Value = $u.F
    $ws.Range("G$r").Value = $u.G
    $ws.Range("H$r").Value = $u.H
}

# Bump the "last updated" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 14:22"
